# This sheet stores its "Price" (column D) values as plain text (e.g. "244.80",
# "0.05770") rather than numbers, so that values with significant trailing
# zeros / fixed decimal formatting render exactly as scraped. Writing a
# numeric-looking string straight into Range.Value/.Value2 makes Excel parse
# it as a real number (losing the trailing zeros and the text type), so for
# each numeric-looking price we:
#   1. temporarily format the cell as Text ("@") so the assignment is kept
#      literally as a string,
#   2. assign the new text,
#   3. restore the cell's style to "Normal" (matches the workbook's original
#      default/unstyled look for these cells).
# Cells whose new value isn't numeric-looking (column E labels) are assigned
# directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.52"
Set-TextValue "D3" "21.78"
Set-TextValue "D4" "5.449"
Set-TextValue "D5" "0.05756"
Set-TextValue "D6" "3.406"
Set-TextValue "D7" "6.311"
Set-TextValue "D8" "0.8153"
Set-TextValue "D9" "1.037"
$ws.Range("E9").Value2 = "8FTXTokenFTTBestin24h"
Set-TextValue "D11" "0.07294"
Set-TextValue "D12" "0.03111"
Set-TextValue "D13" "0.03133"
Set-TextValue "D14" "4.135"
Set-TextValue "D15" "0.09362"
Set-TextValue "D16" "0.001603"
Set-TextValue "D17" "0.04806"
Set-TextValue "D18" "0.0005825"
Set-TextValue "D19" "0.006183"
Set-TextValue "D20" "0.004135"
Set-TextValue "D21" "0.0009914"
Set-TextValue "D23" "3.732"
Set-TextValue "D24" "2.162"
Set-TextValue "D25" "0.3265"
Set-TextValue "D26" "0.1296"
Set-TextValue "D27" "0.0003987"
Set-TextValue "D40" "0.03848"
Set-TextValue "D41" "0.006685"
$ws.Range("E41").Value2 = "40KickTokenKICK"
Set-TextValue "D42" "0.1070"
Set-TextValue "D43" "0.002607"
Set-TextValue "D44" "0.006547"
Set-TextValue "D45" "0.00005583"
Set-TextValue "D46" "0.00000000748"
Set-TextValue "D47" "0.3888"
Set-TextValue "D49" "0.00002094"
Set-TextValue "D50" "0.01007"
